# Add a new "2022-Q1" sheet (same layout as "2021-Q4") positioned between
# "2021-Q4" and "总计", and add a 2022-Q1 summary row (as the new first
# data row) on the "总计" sheet.

$wb = $excel.ActiveWorkbook

$sheetQ4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet right after "2021-Q4" (i.e.
#    right before "总计").
# ---------------------------------------------------------------------
$sheetQ1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$sheetQ1.Name = "2022-Q1"

# NOTE: inserting a sheet shifts the sheet-index-based reference that
# was captured for "总计" before the insert, so it must be re-fetched
# by name afterwards (and again after any further structural change).
$sheetTotal = $wb.Worksheets.Item("总计")

# Copy header row formatting (bold / centered / bordered) from the
# "2021-Q4" sheet so the new sheet matches the existing look & feel.
$sheetQ4.Range("B1:H1").Copy()
$sheetQ1.Range("B1").PasteSpecial(-4122)

# Copy the column-A style (bold / centered / bordered) too.
$sheetQ4.Range("A2:A3").Copy()
$sheetQ1.Range("A2:A3").PasteSpecial(-4122)

# Header values.
$sheetQ1.Range("B1").Value = "基金代码"
$sheetQ1.Range("C1").Value = "基金名称"
$sheetQ1.Range("D1").Value = "基金规模"
$sheetQ1.Range("E1").Value = "股票总仓位"
$sheetQ1.Range("F1").Value = "仓位占比"
$sheetQ1.Range("G1").Value = "持有市值(亿元)"
$sheetQ1.Range("H1").Value = "仓位排名"

# Data rows 2-3: columns B-G hold text-looking numbers (as in the
# source workbook) so force Text format before assigning, then reset
# the cell style back to "no explicit style" (matching row 2 of
# "2021-Q4") via a format-only paste from an already-unstyled cell.
$sheetQ1.Range("B2:G3").NumberFormat = "@"

$sheetQ1.Range("A2").Value = 0
$sheetQ1.Range("B2").Value = "010613"
$sheetQ1.Range("C2").Value = "中融产业趋势一年定期开放混合A"
$sheetQ1.Range("D2").Value = "0.69"
$sheetQ1.Range("E2").Value = "87.17"
$sheetQ1.Range("F2").Value = "5.28"
$sheetQ1.Range("G2").Value = "0.0364"
$sheetQ1.Range("H2").Value = 5

$sheetQ1.Range("A3").Value = 1
$sheetQ1.Range("B3").Value = "010614"
$sheetQ1.Range("C3").Value = "中融产业趋势一年定期开放混合C"
$sheetQ1.Range("D3").Value = "0.11"
$sheetQ1.Range("E3").Value = "87.17"
$sheetQ1.Range("F3").Value = "5.28"
$sheetQ1.Range("G3").Value = "0.0058"
$sheetQ1.Range("H3").Value = 5

# Clean the Text-format artefact off B2:G3 (restore "no explicit
# style", same as the analogous cells on "2021-Q4").
$sheetQ4.Range("B2").Copy()
$sheetQ1.Range("B2:G3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of the data on "总计"
#    (row 2), pushing the existing "2021-Q4" row down to row 3.
# ---------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Rows.Item(2).Insert(-4121)

# Re-apply the column-A style and clear any inherited styling on
# B2:D2 so the row matches the rest of the sheet.
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)

$sheetTotal.Range("B3").Copy()
$sheetTotal.Range("B2:D2").PasteSpecial(-4122)

# Re-number the row index column (A) - row 2 is the new row (index 0),
# row 3 (the shifted former row 2 / "2021-Q4") becomes index 1.
$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("A3").Value = 1

$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 2
$sheetTotal.Range("D2").Value = 0.04
